$d = $word.ActiveDocument

# Helper: find the paragraph whose text starts with a given prefix.
function Get-ParaByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "$prefix*") {
            return $p
        }
    }
    return $null
}

# --- Change 1: append " 1901630" after "MAFERSON DE MELO EVANGELISTA RA:" as
#     its own run. We temporarily park the hidden "_GoBack" bookmark between
#     the two pieces of text (a non-text bookmark marker stops the two runs
#     from being re-coalesced into one run on save); the bookmark is then
#     relocated to its real destination in change 2, so no stray bookmark is
#     left behind. ---
$pMaferson = Get-ParaByPrefix "MAFERSON DE MELO EVANGELISTA RA:"
$rMaferson = $pMaferson.Range.Duplicate
$rMaferson.End = $rMaferson.End - 1
$rMaferson.Collapse(0)
$rMaferson.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rMaferson)
$rMaferson.Text = " 1901630"

# --- Change 2a: move the hidden "_GoBack" bookmark to the end of the
#     "CELULAR: 11 96035-0474" paragraph (right after "MAFERSON ... RA:
#     1901630" / before the empty paragraph + ISRAEL block). ---
$pCelular = Get-ParaByPrefix "CELULAR: 11 96035-0474"
$rCelular = $pCelular.Range.Duplicate
$rCelular.End = $rCelular.End - 1
$rCelular.Collapse(0)
$rCelular.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rCelular)
$rCelular.Text = ""

# --- Change 2b: merge the "PARCERIA:" run and " NIKE DO BRASIL LTDA    " run
#     into a single run (the bookmark that used to sit between them is gone
#     now that it was relocated above). ---
$found = $d.Content.Find.Execute("PARCERIA: NIKE DO BRASIL LTDA    ", $true, $false, $false, $false, $false, $true, 1, $false, "PARCERIA: NIKE DO BRASIL LTDA    ", 2)
